# "Results from R script" — append the latest day's OHLC row pulled in by
# the data refresh, and fix up the previous row's timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 71 (2024-06-03): the date/time value is corrected.
$ws.Cells.Item(71, 1).Value = 45446.2916666667

# New row 72 (2024-06-04): date, volume, high, low, open, close, adj_close, ticker
$ws.Cells.Item(72, 1).Value = 45447.3101967593
$ws.Cells.Item(72, 2).Value = 1500
$ws.Cells.Item(72, 3).Value = 2.97000002861023
$ws.Cells.Item(72, 4).Value = 2.97000002861023
$ws.Cells.Item(72, 5).Value = 2.97000002861023
$ws.Cells.Item(72, 6).Value = 2.97000002861023

# adj_close is stored as text in this sheet (matches the other rows' shared
# string values), so force text formatting before assigning it.
$adjClose = $ws.Cells.Item(72, 7)
$adjClose.NumberFormat = "@"
$adjClose.Value = "2.97000002861023"
$adjClose.Style = "Normal"

$ws.Cells.Item(72, 8).Value = "ESPE.MI"

# Give the new date cell the same date/time display format as the rest of
# column A (copy format only, values already set above).
$ws.Cells.Item(71, 1).Copy()
$ws.Cells.Item(72, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
